$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellRef, $newValue) {
    $rng = $ws.Range($cellRef)
    $rng.NumberFormat = "@"
    $rng.Value = $newValue
    $rng.Style = "Normal"
}

Set-TextValue "D2" "279.47"
Set-TextValue "E2" "0.75%"

Set-TextValue "D3" "27.50"
Set-TextValue "E3" "1.24%"

Set-TextValue "D4" "4.836"
Set-TextValue "E4" "-2.11%"

Set-TextValue "D5" "0.06397"
Set-TextValue "E5" "-0.20%"

Set-TextValue "D6" "7.037"
Set-TextValue "E6" "1.61%"

Set-TextValue "D7" "1.316"
Set-TextValue "E7" "5.22%"

Set-TextValue "D8" "0.9037"
Set-TextValue "E8" "2.53%"

Set-TextValue "D9" "0.1540"
Set-TextValue "E9" "0.94%"

Set-TextValue "D10" "0.06107"
Set-TextValue "E10" "20.93%"

Set-TextValue "D11" "0.07494"
Set-TextValue "E11" "-0.34%"

Set-TextValue "D12" "0.02929"
Set-TextValue "E12" "1.26%"

Set-TextValue "D13" "0.08988"
Set-TextValue "E13" "-0.28%"

Set-TextValue "D14" "0.001577"
Set-TextValue "E14" "0.24%"

Set-TextValue "D15" "0.0006426"
Set-TextValue "E15" "0.42%"

Set-TextValue "D16" "0.006064"
Set-TextValue "E16" "0.24%"

Set-TextValue "D17" "3.492"
Set-TextValue "E17" "1.06%"

Set-TextValue "D18" "3.325"
Set-TextValue "E18" "0.14%"

Set-TextValue "D19" "2.225"
Set-TextValue "E19" "-2.06%"

Set-TextValue "E21" "1.09%"

Set-TextValue "D22" "3.916"
Set-TextValue "E22" "0.17%"

Set-TextValue "D23" "0.04410"
Set-TextValue "E23" "-0.20%"

Set-TextValue "D24" "0.1504"
Set-TextValue "E24" "8.99%"

Set-TextValue "D25" "0.001175"
Set-TextValue "E25" "0.09%"

Set-TextValue "E26" "10.49%"

Set-TextValue "E28" "-1.74%"

Set-TextValue "E29" "-14.43%"

Set-TextValue "D40" "0.04068"
Set-TextValue "E40" "-1.77%"

Set-TextValue "D41" "0.006657"
Set-TextValue "E41" "-2.36%"

Set-TextValue "D42" "0.1394"
Set-TextValue "E42" "18.25%"

Set-TextValue "D43" "0.002090"
Set-TextValue "E43" "-3.31%"

Set-TextValue "D44" "0.01099"
Set-TextValue "E44" "-2.23%"

Set-TextValue "D45" "0.00005546"
Set-TextValue "E45" "6.51%"

Set-TextValue "D46" "1.628"
Set-TextValue "E46" "9.67%"

Set-TextValue "D47" "0.01848"
Set-TextValue "E47" "-8.76%"
